$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the new question, right after "How do you delete a table?"
$ws.Rows("64:64").Insert()
$ws.Range("D64").Value = "How do you delete all data in a table?"
